$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: append the 5 new address strings as rows 4-8 first so the shared
# string table gets them in this specific first-seen order (matches the
# target xlsx's sharedStrings.xml ordering).
$ws.Range("A4").Value = "上海豫园"
$ws.Range("A5").Value = "东方明珠电视塔"
$ws.Range("A6").Value = "上海中心大厦"
$ws.Range("A7").Value = "一大会址"
$ws.Range("A8").Value = "南浦大桥"

# Step 2: rewrite rows 2-8 into the final, reordered sequence.
$ws.Range("A2").Value = "一大会址"
$ws.Range("A3").Value = "上海市人民广场"
$ws.Range("A4").Value = "东方明珠电视塔"
$ws.Range("A5").Value = "南浦大桥"
$ws.Range("A6").Value = "上海中心大厦"
$ws.Range("A7").Value = "上海豫园"
$ws.Range("A8").Value = "上海市正大广场"

# Step 3: extend the "qh" defined name / query-table range to cover the new
# rows but start one row later (A3 instead of A2), matching the diff.
$wb.Names.Item("qh").RefersTo = "=Sheet1!`$A`$3:`$A`$2582"

# Step 4: update the active selection to A7, matching the saved view state.
$ws.Range("A7").Select()
